# Add results from Silverlake.py
#
# The sheet holds several differently-sized pandas index columns
# ("Unnamed: 0", "Unnamed: 0.1", ...) side by side with two SILVER_FOR value
# columns. This change:
#   - inserts two more "Unnamed: 0.x" index columns before the existing ones
#     (shifting the old B..F columns to D..H), and
#   - appends 10 more rows of results to the row-index column and the
#     SILVER_FOR value column that now lives in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns at B:C. Excel shifts the former B,C,D,E,F columns
#    (and all their data/headers) to D,E,F,G,H automatically.
$ws.Range("B1:C1").EntireColumn.Insert()

# 2) Give the two new header cells the same style as the other header cells
#    (bold + boxed), by copying the format from the neighbouring shifted
#    header cell, then set their text.
$ws.Range("D1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("B1").Value = "Unnamed: 0.4"
$ws.Range("C1").Value = "Unnamed: 0.3"

# The column insert also copies column A's boxed style down into the new
# B/C data cells (rows 2-21) because it inherits formatting from the column
# on its left. Strip that back off so those plain data cells carry no
# explicit style, matching the sheet's other (unstyled) data columns.
$ws.Range("B2:C21").ClearFormats()

# 3) Fill the brand-new index columns B (0..24, rows 2-26) and C (0..19, rows
#    2-21) -- the same "0..n-1" pattern used by the pre-existing index
#    columns.
for ($i = 0; $i -le 24; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $i
}
for ($i = 0; $i -le 19; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $i
}

# 4) Extend the plain row-index column A down to row 31 (values 20..29),
#    keeping the same boxed style used by the rest of column A.
$ws.Range("A21").Copy()
$ws.Range("A22:A31").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
for ($i = 0; $i -le 9; $i++) {
    $ws.Cells.Item(22 + $i, 1).Value = 20 + $i
}

# 5) Extend the SILVER_FOR value column (now column H after the insert) with
#    the 10 additional computed values for the new rows.
$newH = @(29.09709106441289, 28.96379896414152, 28.94952882821161, `
    29.52255201468267, 28.89384280964123, 28.95396845664823, `
    28.94541683747048, 28.70188153581586, 27.87425082086361, `
    28.09067066164505)
for ($i = 0; $i -lt $newH.Length; $i++) {
    $ws.Cells.Item(22 + $i, 8).Value = $newH[$i]
}

$ws.Range("A1").Select() | Out-Null
